$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Correção nos dados: os cabeçalhos "unnamed: 1_level_1" e "unnamed: 5_level_1"
# eram rótulos de nível gerados automaticamente pelo pandas; corrigidos para "total".
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
